$wb = $excel.ActiveWorkbook

# Grab reference to the original "ExtTest40mm" sheet.
$orig = $wb.Worksheets.Item("ExtTest40mm")

# Add (and immediately remove) a placeholder sheet so that the internal
# sheetId counter advances past 3; this makes the upcoming Copy() land on
# sheetId 4, matching the id that the new "ExtTest40mm_2" sheet needs to have.
$placeholder = $wb.Worksheets.Add()
$placeholderName = $placeholder.Name

# Duplicate the "ExtTest40mm" sheet; the copy is placed immediately after it.
$orig.Copy($null, $orig)
$copy = $wb.Worksheets.Item(2)

# Remove the placeholder sheet now that the copy has consumed sheetId 4.
# (Re-fetch it by name: inserting the copy shifts indices, so the original
# $placeholder object reference no longer points at the right sheet.)
$wb.Worksheets.Item($placeholderName).Delete()

# Rename the sheets: original becomes "ExtTest40mm_1", the duplicate becomes
# "ExtTest40mm_2".
$orig.Name = "ExtTest40mm_1"
$copy.Name = "ExtTest40mm_2"

# The new "ExtTest40mm_2" sheet represents a not-yet-run test, so clear out
# all of the recorded measurements while leaving the formulas, labels and
# formatting intact.
$copy.Range("C6:P7").ClearContents()
$copy.Range("C8").Clear()
$copy.Range("D8:P8").ClearContents()
$copy.Range("C9:P9").ClearContents()
$copy.Range("C10").ClearContents()
$copy.Range("D10:P10").ClearContents()
$copy.Range("C13:P13").ClearContents()

# Make the new sheet the active tab with C7 selected.
$copy.Activate()
$copy.Range("C7").Select()
